$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(133).Insert()
$ws.Cells.Item(133, 4).NumberFormat = $ws.Cells.Item(134, 4).NumberFormat

$ws.Cells.Item(133, 1).Value = 1
$ws.Cells.Item(133, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(133, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(133, 4).Value = 45180
$ws.Cells.Item(133, 5).Value = 15
$ws.Cells.Item(133, 6).Value = 100112021
$ws.Cells.Item(133, 7).Value = "Ají"
$ws.Cells.Item(133, 8).Value = "Inferno"
$ws.Cells.Item(133, 9).Value = "Primera"
$ws.Cells.Item(133, 10).Value = 250
$ws.Cells.Item(133, 11).Value = 34000
$ws.Cells.Item(133, 12).Value = 35000
$ws.Cells.Item(133, 13).Value = 34600
$ws.Cells.Item(133, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(133, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(133, 16).Value = 2307
$ws.Cells.Item(133, 17).Value = 15
$ws.Cells.Item(133, 18).Value = "Hortaliza"
